$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting rows 18:31 down to 19:32
$ws.Rows.Item(18).Insert()

# Fill the new row 18 with the data for this entry
$ws.Range("A18").Value = 3
$ws.Range("B18").Value = "Femacal de La Calera"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44438
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 100112035
$ws.Range("G18").Value = "Bruselas (repollito)"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 75
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19467
$ws.Range("N18").Value = "$/malla 15 kilos"
$ws.Range("O18").Value = "Provincia de Quillota"
$ws.Range("P18").Value = 1298
$ws.Range("Q18").Value = 15
$ws.Range("R18").Value = "Hortaliza"
